# edit.ps1 -- apply the 'add all the links to fetch texts' change
# Expands column B with the fetched hyperlink texts for each program/club,
# widens column B, wraps the new text, sets matching row heights, removes the
# trailing 'LaunchWeek' row, and moves the lone remaining hyperlink to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Hyperlinks: clear existing (A2, A3, B2), values get rewritten below,
# and only A2 / A3 / B13 should carry a hyperlink afterwards ---
$ws.Hyperlinks.Delete()

# --- Drop the trailing 'LaunchWeek' row (old row 14, had no Link entry) ---
$ws.Rows(14).Delete()

# --- Column B becomes much wider to host the long link lists ---
$ws.Columns("B:B").ColumnWidth = 83

# --- Fill column B (rows 2-13) with the newly fetched link texts ---
$ws.Range("B2").Value = 'https://entrepreneurship.ubc.ca/  | https://entrepreneurship.ubc.ca/engage-us/venture-founder 
| 
https://entrepreneurship.ubc.ca/engage-us/venture-studios
| 
https://entrepreneurship.ubc.ca/engage-us/hatch-venture-builder
|
https://entrepreneurship.ubc.ca/engage-us/internship-program
|
https://entrepreneurship.ubc.ca/engage-us/credit-courses
|
https://entrepreneurship.ubc.ca/engage-us/black-advisory-hub
|
https://entrepreneurship.ubc.ca/engage-us/discovery-foundation
|
https://entrepreneurship.ubc.ca/engage-us/job-board
|
https://entrepreneurship.ubc.ca/engage-us/your-application-journey-0 
| 
https://entrepreneurship.ubc.ca/engage-us/lab2market'
$ws.Range("B3").Value = 'https://entrepreneurship.ok.ubc.ca/ | 
https://entrepreneurship.ok.ubc.ca/programs/ediscovery/ | 
https://entrepreneurship.ok.ubc.ca/programs/startup-sprint/ | 
https://entrepreneurship.ok.ubc.ca/programs/entrepreneur-in-residence/ | 
https://entrepreneurship.ok.ubc.ca/resources/funding-support/ | 
https://entrepreneurship.ok.ubc.ca/resources/employment/ | 
https://entrepreneurship.ok.ubc.ca/resources/equal-opportunity/ | 
https://entrepreneurship.ok.ubc.ca/resources/accelerators/'
$ws.Range("B4").Value = 'http://socialenterpriseclub.ca/ | 
http://socialenterpriseclub.ca/about/ '
$ws.Range("B5").Value = 'https://amscampusbase.ubc.ca/eprojubc/home/'
$ws.Range("B6").Value = 'https://enactusubc.ca/ | 
https://enactusubc.ca/about | 
https://enactusubc.ca/competitions | 
https://enactusubc.ca/ennovate | 
https://enactusubc.ca/enspire | 
https://enactusubc.ca/seeder | 
https://enactusubc.ca/building-bridges '
$ws.Range("B7").Value = 'https://innovation.ubc.ca/about/innovation-hubs | 
https://innovation.ubc.ca/about/graham-lee-innovation-centre | 
https://innovation.ubc.ca/how-engage/innovation-development '
$ws.Range("B8").Value = 'https://summitleaders.ca/ | https://summitleaders.ca/about-us/ | https://summitleaders.ca/programs/ | https://summitleaders.ca/resources/ | https://summitleaders.ca/peer-forward-form/ '
$ws.Range("B9").Value = 'https://www.sauder.ubc.ca/LIFT  | 
https://www.sauder.ubc.ca/current-students/support-involvement-opportunities/ubc-sauder-lift/ubc-sauder-lift-learning-cycle | 
https://www.sauder.ubc.ca/current-students/support-involvement-opportunities/ubc-sauder-lift/our-impact | 
https://www.sauder.ubc.ca/current-students/support-involvement-opportunities/ubc-sauder-lift/our-partners | 
https://www.sauder.ubc.ca/current-students/support-involvement-opportunities/ubc-sauder-lift/student-experience '
$ws.Range("B10").Value = 'https://www.mbasociety.ca/clubs/business-innovation/ | 
https://www.mbasociety.ca/society | 
https://www.mbasociety.ca/ '
$ws.Range("B11").Value = 'https://www.innovationboard.ca/ | 
https://www.innovationboard.ca//learn-more.html | 
https://www.innovationboard.ca//2017/10/04/faq.html '
$ws.Range("B12").Value = 'https://www.ewb.ca/en/ | 
https://www.ewb.ca/en/about-us/our-approach/ | 
https://www.ewb.ca/en/about-us/history/ | 
https://www.ewb.ca/en/what-we-do/investing-in-people/ | 
https://members.ewb.ca/opportunities/caif | 
https://www.ewb.ca/en/what-we-do/investing-in-people/virtual-talent-xchange/ | 
https://www.ewb.ca/en/what-we-do/investing-in-people/junior-fellowship/ | 
https://www.ewb.ca/en/what-we-do/investing-in-people/the-ewb-fellowship/ | 
https://www.ewb.ca/en/what-we-do/investing-in-people/kumvana-fellowship/ | 
https://www.ewb.ca/en/what-we-do/investing-in-people/leadership-development/ | 
https://www.ewb.ca/en/chapters/ | 
https://www.ewb.ca/en/showcase/small-growing-business-ventures/ | 
https://www.ewb.ca/en/showcase/governance-sustainable-services/ | 
https://www.ewb.ca/en/showcase/evolving-engineering/ | 
https://www.ewb.ca/en/what-we-do/investing-in-ventures/ | 
https://www.ewb.ca/en/what-we-do/advocating-for-change/ | 
https://members.ewb.ca/ | 
https://www.ewb.ca/en/showcase/ | 
https://www.ewb.ca/en/what-we-do/advocating-for-change/ '
$ws.Range("B13").Value = 'https://docs.google.com/document/d/1nKYKuUSEwx9d4x-68LYIr451P_VBsC3MhzjVL6ZJCHs/edit '

# --- Wrap text so the multi-line link lists are fully visible ---
$ws.Range("B2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("B6").WrapText = $true
$ws.Range("B7").WrapText = $true
$ws.Range("B8").WrapText = $true
$ws.Range("B9").WrapText = $true
$ws.Range("B10").WrapText = $true
$ws.Range("B11").WrapText = $true
$ws.Range("B12").WrapText = $true
$ws.Range("B13").WrapText = $true

# --- Row heights matching the wrapped content ---
$ws.Rows(2).RowHeight = 323
$ws.Rows(3).RowHeight = 136
$ws.Rows(4).RowHeight = 34
$ws.Rows(6).RowHeight = 119
$ws.Rows(7).RowHeight = 51
$ws.Rows(8).RowHeight = 51
$ws.Rows(9).RowHeight = 153
$ws.Rows(10).RowHeight = 51
$ws.Rows(11).RowHeight = 51
$ws.Rows(12).RowHeight = 323
$ws.Rows(13).RowHeight = 17

# --- Re-create the hyperlinks: e@UBCV / e@UBCO stay on column A,
# the UBC entrepreneurship link now anchors the last row's Google Doc cell ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:e@UBCV")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:e@UBCO")
$ws.Hyperlinks.Add($ws.Range("B13"), "https://entrepreneurship.ubc.ca/")

# Hyperlinks.Add auto-applies the blue/underlined Hyperlink style; A2/A3 stay
# plain in the target sheet (only the font/mail-link behaviour is wanted), so
# restore their original (Normal) look -- B13 is meant to keep the styling.
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Style = "Normal"

# --- View tidy-up: wider window, slightly lower zoom, scrolled near the bottom ---
$ws.Range("B21").Select()
$excel.ActiveWindow.Zoom = 115

Write-Host "edit.ps1 applied"
